# Epolicy Details and Customer Details POM
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Exide" -- turn the 4 key/value rows into a 2-row, 4-col table:
#   row1 = headers (username, password, <blank>, <blank>)
#   row2 = data    (DTDATAENTRYOPR1, exide@1234, <blank>, <blank>)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Exide")

# Drop the old URL / BrowserType rows entirely (also removes the URL hyperlink)
$ws1.Range("A3:B4").Clear()

# Row 2 becomes the data row: username value moves to A2, styled like
# the rest of the data row (no bold); B2 (exide@1234) already holds the
# right value/hyperlink, just needs a matching blank neighbour in C2.
$ws1.Range("A2").Value = "DTDATAENTRYOPR1"
$ws1.Range("A2").Font.Bold = $false
$ws1.Range("C2").Value = ""
$ws1.Range("C2").Font.Underline = $true
$ws1.Range("C2").Font.Color = -65536 + 16711680
$ws1.Range("C2").Font.Color = 16711680

# Row 1 becomes the header row: username/password labels, bold, plus two
# blank-but-bold filler cells to match the new 4-column width
$ws1.Range("B1").Value = "password"
$ws1.Range("A1:D1").Font.Bold = $true

# Clear out row 6's leftover formatting (becomes a fully blank row)
$ws1.Range("A6").Clear()

# Selection / view tweaks
$ws1.Range("C1").Select()

$wb.Worksheets.Item("Exide").Columns("B").ColumnWidth = 11.85546875
$wb.Worksheets.Item("Exide").Columns("C").ColumnWidth = 32.28515625

# ---------------------------------------------------------------------
# Sheet "ProductName" -- widen the 2-col header/data table to 4 columns,
# adding advisorCode / accountNum
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ProductName")

$ws2.Range("A1").Value = "productName"
$ws2.Range("B1").Value = "proposalNo"
$ws2.Range("C1").Value = "advisorCode"
$ws2.Range("D1").Value = "accountNum"
$ws2.Range("A1:D1").Font.Bold = $true

$ws2.Range("C2").Value = 60000086
$ws2.Range("C2").HorizontalAlignment = -4131
$ws2.Range("C2").VerticalAlignment = -4108

$ws2.Range("D2").Value = 1111111111111

$ws2.Range("C6").Select()

$ws2.Columns("C").ColumnWidth = 12.28515625
$ws2.Columns("D").ColumnWidth = 12.140625

# ---------------------------------------------------------------------
# Workbook-level window sizing
# ---------------------------------------------------------------------
$excel.ActiveWindow.Width = 11340
